# "extra code 5/7 p1 and p2"
# Append 3 new student rows (rows 6-8) to the "HocSinh list" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HocSinh list")

# Row 6: id=18, ten="duc thang dep trai", ngay_sinh=37232.708333333336
$ws.Range("A6").Value = 18
$ws.Range("B6").Value = "duc thang dep trai"
$ws.Range("C6").Value = 37232.708333333336

# Row 7: id=25, ten="duc thang bau", ngay_sinh=37135.708333333336
$ws.Range("A7").Value = 25
$ws.Range("B7").Value = "duc thang bau"
$ws.Range("C7").Value = 37135.708333333336

# Row 8: id=26, ten="duc thang bau", ngay_sinh=37135.708333333336
$ws.Range("A8").Value = 26
$ws.Range("B8").Value = "duc thang bau"
$ws.Range("C8").Value = 37135.708333333336

# Copy the date number format (style) from the existing ngay_sinh column
# (C2) onto the newly added cells so they share the same cell style (s="1")
# instead of Excel creating a brand new number format / style entry.
$ws.Range("C2").Copy()
$ws.Range("C6:C8").PasteSpecial(-4122)
